$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.442.89'
$ws.Range("E2").Value = '  +0.69%  '
# Row 3
$ws.Range("D3").Value = '2.250.18'
$ws.Range("E3").Value = '  +0.11%  '
# Row 4
$ws.Range("E4").Value = '  +0.27%  '
# Row 5
$ws.Range("D5").Value = '''307.98'
$ws.Range("E5").Value = '  +0.51%  '
# Row 6
$ws.Range("D6").Value = '''94.61'
$ws.Range("E6").Value = '  -3.09%  '
# Row 7
$ws.Range("D7").Value = '''0.571'
$ws.Range("E7").Value = '  -0.55%  '
# Row 8
$ws.Range("E8").Value = '  +0.19%  '
# Row 9
$ws.Range("E9").Value = '  -0.56%  '
# Row 10
$ws.Range("D10").Value = '''34.94'
$ws.Range("E10").Value = '  -0.91%  '
# Row 11
$ws.Range("E11").Value = '  -0.58%  '
# Row 12
$ws.Range("D12").Value = '''7.23'
$ws.Range("E12").Value = '  -0.37%  '
# Row 13
$ws.Range("E13").Value = '  +0.77%  '
# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.372.23'
$ws.Range("E14").Value = '  +3.81%  '
# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.843'
$ws.Range("E15").Value = '  +1.47%  '
# Row 16
$ws.Range("E16").Value = '  +0.01%  '
# Row 17
$ws.Range("D17").Value = '44.133.26'
$ws.Range("E17").Value = '  +0.40%  '
# Row 18
$ws.Range("D18").Value = '0.0₃0966'
$ws.Range("E18").Value = '  -0.43%  '
# Row 19
$ws.Range("E19").Value = '  -1.03%  '
# Row 20
$ws.Range("E20").Value = '  +2.02%  '
# Row 21
$ws.Range("D21").Value = '''65.96'
$ws.Range("E21").Value = '  +1.26%  '
# Row 22
$ws.Range("D22").Value = '''237.91'
$ws.Range("E22").Value = '  -1.09%  '
# Row 23
$ws.Range("D23").Value = '''3.02'
$ws.Range("E23").Value = '  +2.88%  '
# Row 24
$ws.Range("E24").Value = '  +3.33%  '
# Row 25
$ws.Range("E25").Value = '  -0.16%  '
# Row 26
$ws.Range("D26").Value = '''38.40'
$ws.Range("E26").Value = '  +5.33%  '
# Row 27
$ws.Range("E27").Value = '  +3.79%  '
# Row 28
$ws.Range("E28").Value = '  -1.38%  '
# Row 29
$ws.Range("D29").Value = '''5.98'
$ws.Range("E29").Value = '  -2.56%  '
# Row 30
$ws.Range("E30").Value = '  +0.30%  '
# Row 31
$ws.Range("D31").Value = '''154.77'
$ws.Range("E31").Value = '  -0.91%  '
# Row 32
$ws.Range("D32").Value = '''0.0804'
$ws.Range("E32").Value = '  -1.34%  '
# Row 33
$ws.Range("E33").Value = '  -0.40%  '
# Row 34
$ws.Range("D34").Value = '''3.11'
$ws.Range("E34").Value = '  -10.15%  '
# Row 35
$ws.Range("E35").Value = '  +3.10%  '
# Row 36
$ws.Range("E36").Value = '  +0.78%  '
# Row 37
$ws.Range("E37").Value = '  +0.48%  '
# Row 38
$ws.Range("D38").Value = '''3.48'
$ws.Range("E38").Value = '  +4.19%  '
# Row 39
$ws.Range("D39").Value = '''14.87'
$ws.Range("E39").Value = '  -2.91%  '
# Row 40
$ws.Range("E40").Value = '  +0.27%  '
# Row 41
$ws.Range("E41").Value = '  +0.00%  '
# Row 42
$ws.Range("E42").Value = '  +0.33%  '
# Row 43
$ws.Range("D43").Value = '1.747.68'
$ws.Range("E43").Value = '  -0.60%  '
# Row 44
$ws.Range("E44").Value = '  +1.75%  '
# Row 45
$ws.Range("D45").Value = '''80.85'
$ws.Range("E45").Value = '  -6.49%  '
# Row 46
$ws.Range("B46").Value = 'ordi'
$ws.Range("C46").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D46").Value = '''71.28'
$ws.Range("E46").Value = '  +4.30%  '
# Row 47
$ws.Range("D47").Value = '''100.00'
$ws.Range("E47").Value = '  -0.73%  '
# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '''4.96'
$ws.Range("E48").Value = '  -3.30%  '
# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''1.61'
$ws.Range("E49").Value = '  +6.31%  '
# Row 50
$ws.Range("D50").Value = '''56.08'
$ws.Range("E50").Value = '  +1.90%  '
# Row 51
$ws.Range("D51").Value = '''8.16'
$ws.Range("E51").Value = '  -0.60%  '
